$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add new data rows to Sheet1 (StockProposal test cases)
$ws.Range("A6").Value = "StockProposal"

$ws.Range("B7").Value = "SP_TC_26"
$ws.Range("C7").Value = "It will sync immediately. Do you want to continue?"

$ws.Range("B8").Value = "SP_TC_30"
$ws.Range("D8").Value = "Stock Proposal"

$ws.Range("B9").Value = "SP_TC_31"
$ws.Range("C9").Value = "Saved Successfully."

# Adjust column widths (closest achievable values to the authored widths)
$ws.Columns.Item(1).ColumnWidth = 12.67
$ws.Columns.Item(3).ColumnWidth = 46.67

# Add new Sheet2 after Sheet1 (at the end of the workbook)
$newSheet = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$newSheet.Name = "Sheet2"

# Keep Sheet1 as the active sheet, with C9 selected
$ws.Activate()
$ws.Range("C9").Select() | Out-Null
